$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "51.491.05"
$ws.Range("E2").Value = "  +0.71%  "

Set-TextValue "D3" "2.984.68"
$ws.Range("E3").Value = "  +1.35%  "

Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue "D5" "382.24"
$ws.Range("E5").Value = "  +1.89%  "

Set-TextValue "D6" "103.55"
$ws.Range("E6").Value = "  +2.16%  "

Set-TextValue "D7" "0.547"
$ws.Range("E7").Value = "  +1.79%  "

$ws.Range("E8").Value = "  +0.06%  "

Set-TextValue "D9" "0.592"
$ws.Range("E9").Value = "  +1.12%  "

Set-TextValue "D10" "36.76"
$ws.Range("E10").Value = "  +1.04%  "

$ws.Range("E11").Value = "  -0.85%  "

Set-TextValue "D12" "0.0860"
$ws.Range("E12").Value = "  +1.14%  "

Set-TextValue "D13" "3.461.62"
$ws.Range("E13").Value = "  +1.48%  "

Set-TextValue "D14" "18.44"
$ws.Range("E14").Value = "  +2.06%  "

Set-TextValue "D15" "7.80"
$ws.Range("E15").Value = "  +2.94%  "

Set-TextValue "D16" "2.993.79"
$ws.Range("E16").Value = "  +1.91%  "

Set-TextValue "D17" "11.15"
$ws.Range("E17").Value = "  +4.40%  "

Set-TextValue "D18" "0.998"
$ws.Range("E18").Value = "  +0.28%  "

Set-TextValue "D19" "51.512.08"
$ws.Range("E19").Value = "  +0.86%  "

Set-TextValue "D20" "3.08"
$ws.Range("E20").Value = "  -0.62%  "

Set-TextValue "D21" "12.60"
$ws.Range("E21").Value = "  +1.24%  "

Set-TextValue "D22" "0.0₃0962"
$ws.Range("E22").Value = "  +0.46%  "

Set-TextValue "D23" "70.55"
$ws.Range("E23").Value = "  +2.67%  "

Set-TextValue "D24" "267.71"
$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("E25").Value = "  +2.00%  "

Set-TextValue "D26" "7.85"
$ws.Range("E26").Value = "  -4.11%  "

Set-TextValue "D27" "7.45"
$ws.Range("E27").Value = "  -2.11%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  +1.82%  "

Set-TextValue "D30" "26.03"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("E31").Value = "  -0.55%  "

Set-TextValue "D32" "10.31"
$ws.Range("E32").Value = "  +2.88%  "

Set-TextValue "D33" "34.68"
$ws.Range("E33").Value = "  +3.87%  "

Set-TextValue "D34" "51.55"
$ws.Range("E34").Value = "  +1.64%  "

Set-TextValue "D35" "2.07"
$ws.Range("E35").Value = "  +1.09%  "

$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("E38").Value = "  +2.91%  "

Set-TextValue "D39" "16.81"
$ws.Range("E39").Value = "  +2.52%  "

Set-TextValue "D40" "0.116"
$ws.Range("E40").Value = "  +1.35%  "

$ws.Range("E41").Value = "  +2.75%  "

Set-TextValue "D42" "2.55"
$ws.Range("E42").Value = "  +2.84%  "

Set-TextValue "D43" "124.78"
$ws.Range("E43").Value = "  +3.58%  "

Set-TextValue "D44" "3.68"
$ws.Range("E44").Value = "  +10.22%  "

Set-TextValue "D45" "21.46"
$ws.Range("E45").Value = "  +0.69%  "

# Rows 46 and 47 swap content (ApeXProtocol <-> WEMIXToken)
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D46" "2.03"
$ws.Range("E46").Value = "  +0.15%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D47" "2.38"
$ws.Range("E47").Value = "  +3.44%  "

$ws.Range("E48").Value = "  -0.95%  "

Set-TextValue "D49" "2.041.64"
$ws.Range("E49").Value = "  +2.56%  "

$ws.Range("E50").Value = "  +16.23%  "

Set-TextValue "D51" "0.0331"
$ws.Range("E51").Value = "  +2.06%  "
